$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows for years 2000-2009 (original rows 2-11).
# After this, the old 2010-2020 rows (12-22) shift up to become rows 2-12.
$ws.Range("A2:F11").EntireRow.Delete() | Out-Null

# Clean up a handful of values that had excess floating point precision
# in the 2010-2016 and 2020 rows (now rows 2-8 and 12).
$ws.Cells.Item(2, 4).Value2 = 65343        # D2 2010 female population
$ws.Cells.Item(2, 6).Value2 = 68748        # F2 2010 male population
$ws.Cells.Item(3, 5).Value2 = 134916       # E3 2011 total population
$ws.Cells.Item(4, 5).Value2 = 135922       # E4 2012 total population
$ws.Cells.Item(5, 5).Value2 = 136726       # E5 2013 total population
$ws.Cells.Item(6, 5).Value2 = 137646       # E6 2014 total population
$ws.Cells.Item(7, 5).Value2 = 138326       # E7 2015 total population
$ws.Cells.Item(8, 5).Value2 = 139232       # E8 2016 total population
$ws.Cells.Item(12, 2).Value2 = 50992       # B12 2020 rural population
$ws.Cells.Item(12, 3).Value2 = 90220       # C12 2020 urban population

# Add the new 2021 and 2022 rows (13 and 14), copying the formatting of
# the preceding data row so the year label keeps its style.
$ws.Range("A12:F12").Copy() | Out-Null
$ws.Range("A13:F14").PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(13, 1).Value2 = "2021年"
$ws.Cells.Item(13, 2).Value2 = 49835
$ws.Cells.Item(13, 3).Value2 = 91425
$ws.Cells.Item(13, 4).Value2 = 68949
$ws.Cells.Item(13, 5).Value2 = 141260
$ws.Cells.Item(13, 6).Value2 = 72311

$ws.Cells.Item(14, 1).Value2 = "2022年"
$ws.Cells.Item(14, 2).Value2 = 49104
$ws.Cells.Item(14, 3).Value2 = 92071
$ws.Cells.Item(14, 4).Value2 = 68969
$ws.Cells.Item(14, 5).Value2 = 141175
$ws.Cells.Item(14, 6).Value2 = 72206

$excel.CutCopyMode = 0
